# "Chrome and firefox issue sorted"
#
# Adds a new "LoginPageUI" worksheet (holding the UI text used to verify the
# login page renders correctly across browsers) positioned right before the
# existing "InvalidPasswordLogin" sheet, and leaves the selection on the
# new sheet / repoints a couple of selections that a real Excel session
# would naturally update when the sheet set changes.

$wb = $excel.ActiveWorkbook

# --- loginData sheet: drop the now-unused extra formatted column and move
#     the active selection to A2 (no longer the active/selected tab).
$loginData = $wb.Worksheets.Item(1)
[void]$loginData.Columns("G").Delete()
[void]$loginData.Range("A2").Select()

# --- Insert the new "LoginPageUI" worksheet immediately before
#     "InvalidPasswordLogin" (keeps loginData / InvalidUserNameLogin in
#     place, and InvalidPasswordLogin stays last).
$invalidPasswordLogin = $wb.Worksheets.Item("InvalidPasswordLogin")
$loginPageUI = $wb.Worksheets.Add($invalidPasswordLogin)
$loginPageUI.Name = "LoginPageUI"

# Populate the header row and the sample data row. Values are entered in
# this order so that new shared-string entries line up with the authored
# workbook's string table order.
$loginPageUI.Range("B2").Value = "Log on using your details"
$loginPageUI.Range("C2").Value = "View Terms and Conditions"
$loginPageUI.Range("A1").Value = "loginpagetitle"
$loginPageUI.Range("B1").Value = "headerlogonmetext"
$loginPageUI.Range("C1").Value = "verifytermsandconditions"
$loginPageUI.Range("A2").Value = "Log On - Ci Anywhere"

# Leave this new sheet active, with the same stray selection the authored
# workbook shows.
[void]$loginPageUI.Range("E22").Select()
